$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 359-360; this shifts the existing rows 359..415
# down to 361..417 (and the dimension/ref grows to T417 automatically).
$ws.Rows("359:360").Insert()

# --- New row 359 ---
$ws.Cells.Item(359, 1).Value = 1
$ws.Cells.Item(359, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(359, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(359, 4).Value = 45142
$ws.Cells.Item(359, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(359, 5).Value = 15
$ws.Cells.Item(359, 6).Value = "Fruta"
$ws.Cells.Item(359, 7).Value = 100108
$ws.Cells.Item(359, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(359, 9).Value = 100108006
$ws.Cells.Item(359, 10).Value = "Plátano"
$ws.Cells.Item(359, 11).Value = "Sin especificar"
$ws.Cells.Item(359, 12).Value = "Pintón"
$ws.Cells.Item(359, 13).Value = 150
$ws.Cells.Item(359, 14).Value = 17000
$ws.Cells.Item(359, 15).Value = 17000
$ws.Cells.Item(359, 16).Value = 17000
$ws.Cells.Item(359, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(359, 18).Value = "Ecuador"
$ws.Cells.Item(359, 19).Value = 850
$ws.Cells.Item(359, 20).Value = 20

# --- New row 360 ---
$ws.Cells.Item(360, 1).Value = 1
$ws.Cells.Item(360, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(360, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(360, 4).Value = 45142
$ws.Cells.Item(360, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(360, 5).Value = 15
$ws.Cells.Item(360, 6).Value = "Fruta"
$ws.Cells.Item(360, 7).Value = 100108
$ws.Cells.Item(360, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(360, 9).Value = 100108006
$ws.Cells.Item(360, 10).Value = "Plátano"
$ws.Cells.Item(360, 11).Value = "Sin especificar"
$ws.Cells.Item(360, 12).Value = "Verde"
$ws.Cells.Item(360, 13).Value = 150
$ws.Cells.Item(360, 14).Value = 18000
$ws.Cells.Item(360, 15).Value = 18000
$ws.Cells.Item(360, 16).Value = 18000
$ws.Cells.Item(360, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(360, 18).Value = "Ecuador"
$ws.Cells.Item(360, 19).Value = 900
$ws.Cells.Item(360, 20).Value = 20
